$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the table to accommodate the new rows (A1:B2 -> A1:B5)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B5"))

# Existing row 2: update the folder reference to something more realistic
$ws.Range("B2").Value = "10-tourisme"

# New row 3: bevnat_info belongs to the bevnat folder
$ws.Range("B3").Value = "bevnat"
$ws.Range("A3").Value = "bevnat_info"

# New row 4: bevnat_variable also belongs to the bevnat folder
$ws.Range("B4").Value = "bevnat"
$ws.Range("A4").Value = "bevnat_variable"

# New row 5: statpop_info belongs to the statpop folder
$ws.Range("B5").Value = "statpop"
$ws.Range("A5").Value = "statpop_info"

# Move the selection below the new data, matching the saved view state
$ws.Range("A6").Select()
